# Updates the cryptos list sheet to the latest scraped snapshot.
# Values are written as text (a leading apostrophe is used for entries that
# look numeric, e.g. "436.50", so Excel keeps them as text instead of
# auto-converting them to a Double and losing the trailing zero / exact digits).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.955.89"
$ws.Range("E2").Value = "  -4.52%  "
$ws.Range("D3").Value = "3.026.92"
$ws.Range("E3").Value = "  -5.72%  "
$ws.Range("D5").Value = "'579.81"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "'128.42"
$ws.Range("E6").Value = "  -6.56%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.025.55"
$ws.Range("E8").Value = "  -5.54%  "
$ws.Range("D9").Value = "'0.501"
$ws.Range("E9").Value = "  -2.85%  "
$ws.Range("E10").Value = "  -6.58%  "
$ws.Range("D11").Value = "'5.20"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "'0.444"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "'0.0000227"
$ws.Range("E13").Value = "  -6.27%  "
$ws.Range("D14").Value = "'32.99"
$ws.Range("E14").Value = "  -5.73%  "
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "3.520.00"
$ws.Range("E16").Value = "  -5.86%  "
$ws.Range("D17").Value = "3.016.34"
$ws.Range("E17").Value = "  -5.65%  "
$ws.Range("D18").Value = "60.689.56"
$ws.Range("E18").Value = "  -4.66%  "
$ws.Range("D19").Value = "'6.48"
$ws.Range("E19").Value = "  -2.30%  "
$ws.Range("D20").Value = "'436.50"
$ws.Range("E20").Value = "  -6.66%  "
$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = "  -5.80%  "
$ws.Range("D22").Value = "'0.669"
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("D23").Value = "'7.09"
$ws.Range("E23").Value = "  -8.01%  "
$ws.Range("D24").Value = "'12.96"
$ws.Range("E24").Value = "  -4.65%  "
$ws.Range("D25").Value = "'79.86"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").Value = "'7.42"
$ws.Range("E29").Value = "  -5.57%  "
$ws.Range("D30").Value = "'1.93"
$ws.Range("E30").Value = "  -6.65%  "
$ws.Range("D31").Value = "'6.28"
$ws.Range("E31").Value = "  -8.77%  "
$ws.Range("D32").Value = "'25.50"
$ws.Range("E32").Value = "  -7.54%  "
$ws.Range("D33").Value = "'0.0950"
$ws.Range("E33").Value = "  -8.45%  "
$ws.Range("D34").Value = "'2.19"
$ws.Range("E34").Value = "  -9.58%  "
$ws.Range("D35").Value = "'0.962"
$ws.Range("E35").Value = "  -7.66%  "
$ws.Range("D36").Value = "'5.67"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("D37").Value = "'50.20"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("D38").Value = "0.0₃0680"
$ws.Range("E38").Value = "  -7.05%  "
$ws.Range("D39").Value = "'8.56"
$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("D40").Value = "'0.0364"
$ws.Range("E40").Value = "  -7.14%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("D42").Value = "'387.71"
$ws.Range("E42").Value = "  -5.25%  "
$ws.Range("D43").Value = "'2.55"
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("D44").Value = "2.671.43"
$ws.Range("E44").Value = "  -5.66%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "'0.239"
$ws.Range("E47").Value = "  -7.17%  "
$ws.Range("D48").Value = "'118.49"
$ws.Range("E48").Value = "  -6.71%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'24.19"
$ws.Range("E49").Value = "  -6.02%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.108"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("E51").Value = "  +4.34%  "
